# Update "想去人数" (F column) attendance counts on the two data sheets
# that contain the full row set (展览 / sheet1 and 全部类型 / sheet4),
# reflecting a refreshed data scrape (gh-pages output regenerated).

$wb = $excel.ActiveWorkbook

$rowUpdates = @{
    2  = 1050
    3  = 741
    4  = 257
    5  = 30
    8  = 1675
    9  = 6145
    11 = 361
    12 = 289
    13 = 89
    14 = 364
    15 = 134
    16 = 5480
    17 = 267
    18 = 1272
    19 = 137
    20 = 113
    22 = 102
    23 = 264
    24 = 98
    26 = 8
    28 = 2
    29 = 387
    30 = 79
    31 = 51
    32 = 76
    34 = 61
    35 = 22
    36 = 62
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($rowNum in $rowUpdates.Keys) {
        $ws.Cells.Item($rowNum, 6).Value = $rowUpdates[$rowNum]
    }
}
